$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Personal/Carers Leave" row (row 2); remaining rows shift up.
$ws.Rows.Item(2).Delete()

# The deleted row shrank the used range referenced by the HTML_* defined
# names by one row.
$wb.Names.Item("HTML_1").RefersTo = "='Digio Sheet'!`$A`$1:`$F`$4"
$wb.Names.Item("HTML_all").RefersTo = "='Digio Sheet'!`$A`$1:`$F`$4"
